# Update column F (dSF) values on Sheet1 to reflect repulled / recalculated data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = 1
    3  = 1
    6  = -3
    11 = -4
    13 = -7
    15 = 0
    18 = 3
    20 = 1
    24 = 1
    27 = 0
    31 = 0
    32 = 6
    33 = -2
    35 = -4
    37 = 0
    38 = 1
    40 = -2
    43 = 3
    49 = -1
    50 = 3
    58 = 0
    59 = 5
    60 = 1
    61 = 1
    62 = -2
    63 = 1
    64 = -4
    65 = 3
    66 = -6
    67 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
